# Game Boy Printer Emulator - Games Support
# Add two newly-documented games ("Pokémon Picross (Japan, unreleased)" and
# "Hello Kitty Pocket Camera (Japan, unreleasd)") to the games-support table.
# They are inserted as new rows 41 and 42 (pushing the existing rows 41+
# down by two), matching the author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows above the current row 41 so the existing data
# (previously starting at row 41) shifts down to rows 43+.
$ws.Rows.Item(41).Insert()
$ws.Rows.Item(41).Insert()

# Row 41: Pokémon Picross (Japan, unreleased)
$ws.Cells.Item(41, 1).Value = "Pokémon Picross (Japan, unreleased)"
$ws.Cells.Item(41, 2).Value = "GBC"
$ws.Cells.Item(41, 3).Value = "Yes"
$ws.Cells.Item(41, 4).Value = "Normal"
$ws.Cells.Item(41, 5).Value = "No"
$ws.Cells.Item(41, 6).Value = "Standard"

# Row 42: Hello Kitty Pocket Camera (Japan, unreleasd)
$ws.Cells.Item(42, 1).Value = "Hello Kitty Pocket Camera (Japan, unreleasd)"
$ws.Cells.Item(42, 2).Value = "GBC"
$ws.Cells.Item(42, 3).Value = "Yes"
$ws.Cells.Item(42, 4).Value = "Normal"
$ws.Cells.Item(42, 5).Value = "No"
$ws.Cells.Item(42, 6).Value = "Custom"

# Match the style used for the new rows' title cells (same "highlighted"
# look used further down the sheet for other not-yet-confirmed entries).
$ws.Cells.Item(41, 1).Style = $ws.Cells.Item(53, 1).Style
$ws.Cells.Item(42, 1).Style = $ws.Cells.Item(53, 1).Style

# Restore the cursor position left by the author after the edit.
$ws.Range("A63").Select()
